$wb = $excel.ActiveWorkbook

# --- Service Contacts: move "delivery_organisation_path" column so it sits
# immediately before "practitioner_key" (was the last data column, R) ---
$wsSC = $wb.Worksheets.Item("Service Contacts")

# Insert a blank column at D (practitioner_key and everything after shifts right).
$wsSC.Columns("D").Insert()

# After the insert, the old "delivery_organisation_path" column (R) is now at S.
# Move its values (and their formatting) into the newly created column D.
$wsSC.Range("S1:S3").Cut($wsSC.Range("D1:D3"))

# Remove the now-empty source column so everything after it shifts back left.
$wsSC.Columns("S").Delete()

# The column width (19 chars) followed the old column; restore it on D.
$wsSC.Columns("D").ColumnWidth = 18.17

# Update the view/selection state on this sheet.
$wsSC.Columns("D").Select()

# --- View/selection-only changes on other sheets ---
$wsK10 = $wb.Worksheets.Item("K10+")
$wsK10.Range("F1:F5").Select()

$wsK5 = $wb.Worksheets.Item("K5")
$wsK5.Range("F1:F5").Select()

# Organisations is selected last so it remains the active tab, matching the
# original workbook state.
$wsOrg = $wb.Worksheets.Item("Organisations")
$wsOrg.Range("H1:J3").Select()
